$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of score-range / grade-level data. Writing column by
# column (A14, A15, then B14, B15) makes the new entries land in the shared
# string table in that same order.
$ws.Range("A14").Value = "500~550"
$ws.Range("A15").Value = "550~600"
$ws.Range("B14").Value = "대학생 수준"
$ws.Range("B15").Value = "성인 고급 수준"

# Copy the formatting (font/alignment/border) from the last existing table
# row (A13:B13) onto the two new rows so they keep the same look-and-feel.
$src = $ws.Range("A13:B13")
$dst = $ws.Range("A14:B15")
$src.Copy()
$dst.PasteSpecial(-4122)

# These new rows now close off the table, so drop the top/bottom edges on
# one of the cells, leaving just the left/right divider between the two
# columns ...
$a14 = $ws.Range("A14")
$a14.Borders.Item(8).LineStyle = -4142
$a14.Borders.Item(9).LineStyle = -4142

# ... then fan that exact formatting out to the other three new cells so
# they all share one consistent style.
$a14.Copy()
foreach ($addr in @("B14", "A15", "B15")) {
  $ws.Range($addr).PasteSpecial(-4122)
}

# Move the active selection to the new last cell.
$ws.Range("B15").Select()
